$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Students")

# Temporarily force Text number format so the numeric-/date-looking
# strings below are stored as literal text (matching the original
# shared-string cell type) instead of being re-interpreted as
# numbers/dates; ClearFormats afterwards drops the cell-level style
# again so the written cells keep the default style, like the source.
$ws.Range("A2:M2").NumberFormat = "@"
$ws.Range("A2").Value = "1"
$ws.Range("B2").Value = "2025-01-01"
$ws.Range("C2").Value = "satyam"
$ws.Range("D2").Value = "testfather"
$ws.Range("E2").Value = "ramkrishnanagar"
$ws.Range("F2").Value = "7250585058"
$ws.Range("G2").Value = "06:00-10:00, 22:00-06:00"
$ws.Range("H2").Value = "1,2"
$ws.Range("I2").Value = "1"
$ws.Range("J2").Value = "350.00"
$ws.Range("K2").Value = "150.00"
$ws.Range("L2").Value = "1"
$ws.Range("M2").Value = "2025-03-05"
$ws.Range("A2:M2").ClearFormats()

# Remove the now-obsolete data rows (former rows 3 and 4)
$ws.Rows.Item(3).Delete()
$ws.Rows.Item(3).Delete()
